$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E16:E51) is re-entered in descending order
# (it used to run ascending 1610..1909; now it runs descending 1909..1610).
# This mirrors the author re-keying/re-pasting the period list, which in
# turn reorders the shared-strings table for these labels.
$periods = @(
  "1909","1908","1907","1906","1905","1904","1903","1902","1901",
  "1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801",
  "1712","1711","1710","1709","1708","1707","1706","1705","1704","1703","1702","1701",
  "1612","1611","1610"
)

for ($i = 0; $i -lt $periods.Length; $i++) {
  $row = 16 + $i
  $ws.Range("E$row").Value = $periods[$i]
}

# The "Valor Mora" figure tied to period 1909 (24534) and the figure tied
# to period 1610 (32000) swap which row they sit on following the
# re-ordering above.
$ws.Range("F16").Value = 24534
$ws.Range("F51").Value = 32000
